# Updates the "grid_cell" lookup column (AG) on the "solar" worksheet so
# that each row of the distribution table (rows 4-28) is tagged with the
# grid cell that the regenerated model data now associates it with.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$newGridCells = @(
    "CHE_0",
    "CHE_17",
    "CHE_19",
    "CHE_23",
    "CHE_9",
    "CHE_21",
    "CHE_4",
    "CHE_24",
    "CHE_8",
    "CHE_5",
    "CHE_7",
    "CHE_10",
    "CHE_22",
    "CHE_12",
    "CHE_3",
    "CHE_2",
    "CHE_14",
    "CHE_18",
    "CHE_20",
    "CHE_1",
    "CHE_6",
    "CHE_11",
    "CHE_15",
    "CHE_25",
    "CHE_13"
)

$startRow = 4
for ($i = 0; $i -lt $newGridCells.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 33).Value = $newGridCells[$i]
}
